# Adds a second copy of the small redaction picture ("Kép 2") on slide 8
# and nudges the original copy's position, per the target diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# Locate the existing "Kép 2" picture shape on the slide.
$orig = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Kép 2") {
        $orig = $candidate
    }
}

# Duplicate it in-place (this keeps the same embedded image) and then
# reposition both the original and the new copy to their target spots.
$copy = $orig.Duplicate()

# Left/Top/Width/Height round-trip through a single-precision float in the
# PowerPoint object model, so the literals below are the exact float32
# values that convert back to the exact target EMUs (avoids +/-1 EMU drift
# from naive EMU/12700 division).
$orig.Left = 237.84197998046875
$orig.Top = 231.53701782226562
$orig.Width = 51.00708770751953
$orig.Height = 23.25322914123535

$copy.Left = 146.20559692382812
$copy.Top = 148.98709106445312
$copy.Width = 40.99441146850586
$copy.Height = 18.688661575317383
